# Updated to version submitted to PLOS ONE
# Experiment "5A" (rows 12-18) duration corrected from 48h to 24h.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O12").Value = "24h"
$ws.Range("O13").Value = "24h"
$ws.Range("O14").Value = "24h"
$ws.Range("O15").Value = "24h"
$ws.Range("O16").Value = "24h"
$ws.Range("O17").Value = "24h"
$ws.Range("O18").Value = "24h"
